$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.109.32"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.007.23"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.04"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.82"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.006.57"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.23"
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.34"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "3.581.74"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "62.109.93"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "2.993.23"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "446.65"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.13"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.04"
$ws.Range("E25").Value = "  +10.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.37"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.82"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.11"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.47"
$ws.Range("E43").Value = "  +10.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.281"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.40"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "2.717.50"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.74"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -1.59%  "
